$d = $word.ActiveDocument

# First paragraph edit: replace the tail about "so Government is planning..."
# with the new sentence about Link NYC - Citybridge.
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Text = "so Government is planning to Install additional WIFI hotspot across the city."
$find1.Replacement.ClearFormatting()
$find1.Replacement.Text = "Link NYC " + [char]0x2013 + " Citybridge is one of the popular internet service provider has been asked to resolve the issue by installing more hotspot in the NYC."
$find1.Execute([ref]$find1.Text, $true, $false, $false, $false, $false, $true, 1, $false, [ref]$find1.Replacement.Text, 2) | Out-Null

# Second paragraph edit: "official" -> "Service Provider"
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Text = "official wanted to know"
$find2.Replacement.ClearFormatting()
$find2.Replacement.Text = "Service Provider wanted to know"
$find2.Execute([ref]$find2.Text, $true, $false, $false, $false, $false, $true, 1, $false, [ref]$find2.Replacement.Text, 2) | Out-Null

# Second paragraph edit: "Internet Usage." -> "Internet Usage so that Link NYC - Citybridge can use the data for the installation purpose."
$find3 = $d.Content.Find
$find3.ClearFormatting()
$find3.Text = "Internet Usage."
$find3.Replacement.ClearFormatting()
$find3.Replacement.Text = "Internet Usage so that Link NYC " + [char]0x2013 + " Citybridge can use the data for the installation purpose."
$find3.Execute([ref]$find3.Text, $true, $false, $false, $false, $false, $true, 1, $false, [ref]$find3.Replacement.Text, 2) | Out-Null
